# Add a new "2020" column (K) to the table, mirroring the formatting of the
# existing 2019 column (J) and filling in the new values for 2020.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column J (rows 3-8, which carries the border/number styles used by the
# whole data block) into the new column K so the new cells inherit the same
# styles as their row/column neighbours.
$ws.Range("J3:J8").Copy($ws.Range("K3:K8"))

# Row 5 ("Итого") has no value in the new 2020 column, so clear the cell that
# the copy above created there.
$ws.Range("K5").Clear()

# Fill in the actual 2020 figures.
$ws.Range("K4").Value = 2020
$ws.Range("K6").Value = 5.9
$ws.Range("K7").Value = 1.5
$ws.Range("K8").Value = "-"

# Match the saved selection state recorded in the workbook.
$ws.Range("L16").Select()
